$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C, rows 2 through 140 all hold the same serial date value 45171
# which needs to be updated to 45172 (one day later).
for ($r = 2; $r -le 140; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}
